$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("On-Site")
$ws2 = $wb.Worksheets.Item("Replacement")

# 1) Copy formatting (styles) from On-Site rows 2:3 onto Replacement rows 2:3.
#    Done in two pieces (A:T, then V) because column U has no cell on the
#    source row (it's skipped), so a single A:V copy would wrongly stamp U.
$ws1.Range("A2:T3").Copy()
$ws2.Range("A2").PasteSpecial(-4122)
$ws1.Range("V2:V3").Copy()
$ws2.Range("V2").PasteSpecial(-4122)

# 2) Copy the values (preserves text/number/bool typing exactly, e.g. the
#    literal text "False" in column D must stay a shared string, not become
#    a real boolean).
$ws1.Range("A2:T3").Copy()
$ws2.Range("A2").PasteSpecial(-4163)
$ws1.Range("V2:V3").Copy()
$ws2.Range("V2").PasteSpecial(-4163)

$excel.CutCopyMode = $false
